# Quick-quote / sign-up credential refresh on the "CreateAccount" sheet.
# Replaces the previously-recorded Selenium-generated e-mail / password
# values in E2:E6 and F7 with a newer batch, and tightens column E's
# width to fit the new (slightly shorter) longest value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateAccount")

$ws.Range("E2").Value = "SeleniumCjti@mailinator.com"
$ws.Range("E3").Value = "SeleniumPGBr@mailinator.com"
$ws.Range("E4").Value = "SeleniumctfD@mailinator.com"
$ws.Range("E5").Value = "SeleniumSWbf@mailinator.com"
$ws.Range("E6").Value = "SeleniumPfEJ@mailinator.com"
$ws.Range("F7").Value = "Automation6223!"

$ws.Columns.Item(5).ColumnWidth = 28.833
